# Edit the "area_pop_sum" sheet to reshape the summary table:
#   - remove the "Density"/C column entirely
#   - lower-case "Population" -> "population"
#   - add a new 4th row "density" with the value that used to be in C2/C3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area_pop_sum")

# Capture the density value (currently duplicated in C2/C3) before we clear it.
$densityValue = $ws.Range("C2").Value2

# Rename "Population" label to lowercase "population"
$ws.Range("A3").Value2 = "population"

# Add the new 4th row holding the density label/value.
$ws.Range("A4").Value2 = "density"
$ws.Range("B4").Value2 = $densityValue

# Drop the now-unused Density column (C1:C3).
$ws.Range("C1:C3").Clear()
